# Lab2: populate the "flops" worksheet with data and add a scatter chart
# plotting FLOPS vs. number of processors, mirroring the "bandwidth" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("flops")

# --- Headers -----------------------------------------------------------
$ws.Range("A1").Value = "PROCESSORS"
$ws.Range("B1").Value = "FLOPS"

# --- Data ----------------------------------------------------------------
$procs = @(1, 2, 3, 4, 5, 6, 7, 8)
$flops = @(
    671625940.75260198,
    654235532.67820895,
    1397635454.8483801,
    2094010983.5247099,
    2821597040.0268998,
    3496710295.9566398,
    4158952900.3470502,
    4865781902.5522003
)

for ($i = 0; $i -lt $procs.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $procs[$i]
    $ws.Cells.Item($row, 2).Value = $flops[$i]
}

# --- Chart -----------------------------------------------------------
$chartObj = $ws.ChartObjects().Add(368300, 165100, 4445000, 1143000)
$chart = $chartObj.Chart
$chart.ChartType = -4169 # xlXYScatterLines

$ser = $chart.SeriesCollection.NewSeries()
$ser.XValues = $ws.Range("A2:A9")
$ser.Values = $ws.Range("B2:B9")

$chart.HasTitle = $true
$chart.ChartTitle.Text = "Floating Point Operations per Second vs. Processors employed"

$chart.Axes(1).HasTitle = $true
$chart.Axes(1).AxisTitle.Text = "Number of Processors"

$chart.Axes(2).HasTitle = $true
$chart.Axes(2).AxisTitle.Text = "FLOPS"

# --- Active tab / selection -----------------------------------------------
$ws.Range("A2:B9").Select()
$excel.ActiveWindow.ActiveSheet.Activate()
$wb.Worksheets.Item("flops").Activate()
